# Roll the "Chart" sheet's 90-day window forward by one day:
#   - drop 2025-10-08, shift every remaining date back one row, append 2026-01-05
#   - the "Valid" (C) counts follow the same row shift (row 90 keeps its value,
#     since there is no newer count to pull in)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$firstRow = 2
$lastRow = 90

# Capture the current (pre-edit) Date (A) and Valid (C) columns.
$oldDates = @{}
$oldCounts = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $oldDates[$r] = [string]$ws.Cells.Item($r, 1).Value2
    $oldCounts[$r] = $ws.Cells.Item($r, 3).Value2
}

# Make sure the date column keeps storing plain text (otherwise Excel will
# helpfully reinterpret "2025-10-09" as a date serial number).
$ws.Range("A$firstRow`:A$lastRow").NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $d = [datetime]::ParseExact($oldDates[$r], "yyyy-MM-dd", $null)
    $ws.Cells.Item($r, 1).Value = $d.AddDays(1).ToString("yyyy-MM-dd")

    if ($r -lt $lastRow) {
        $ws.Cells.Item($r, 3).Value = $oldCounts[$r + 1]
    }
}
